$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a new "7 FA" results block (rows 66-77), mirroring the layout used by
# the existing "2 FA".."6 FA" blocks above it (header row + 10 metric rows),
# but left empty (no computed values yet) - matching the "5 FA" block's
# formatting (rows 40-51), which is the template this new block is based on.
# ---------------------------------------------------------------------------

# Copy the cell formatting from the "5 FA" block (A40:D51) down onto the new
# "7 FA" block (A66:D77).
$ws.Range("A40:D51").Copy()
$ws.Range("A66:D77").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Block header label.
$ws.Range("A66").Value = "7 FA"

# Column headers for the metrics table.
$ws.Range("A67").Value = "fm"
$ws.Range("B67").Value = "minchi"
$ws.Range("C67").Value = "minrank"
$ws.Range("D67").Value = "ml"

# Row labels (values for B:D are intentionally left blank - data not yet
# collected for this factor count).
$ws.Range("A68").Value = "BIC"
$ws.Range("A69").Value = "TLI"
$ws.Range("A70").Value = "rms "
$ws.Range("A71").Value = "STATISTIC"
$ws.Range("A72").Value = "PVAL"
$ws.Range("A73").Value = "objective"
$ws.Range("A74").Value = "EBIC "
$ws.Range("A75").Value = "dof"
$ws.Range("A76").Value = "chi"
$ws.Range("A77").Value = "RMSEA"

# ---------------------------------------------------------------------------
# Update the view so the newly added block is the one in focus, like it was
# right after the block was entered.
# ---------------------------------------------------------------------------
$ws.Range("E74").Select()
$win = $wb.Windows.Item(1)
$win.ScrollRow = 64
$win.ScrollColumn = 1
